$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the user id embedded both in the free-text note (F2) and the
# dedicated USERID cell (G2).
$ws.Range("F2").Value = "Username : 31160;`nPassword : bni1234;`nTgl. Market : 24/01/2023;`nFile Excel : 24012023HargaPasarFixedIncome.xlsx"
$ws.Range("G2").Value = 31160

# O2 keeps the same text value (file name) - no change needed, but set it
# explicitly to be safe / idempotent.
$ws.Range("O2").Value = "24012023HargaPasarFixedIncome.xlsx"

# Update the active selection/view: no more frozen/scrolled topLeftCell,
# and selection moves from G3 to F3.
$ws.Range("A1").Select()
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("F3").Select()
